$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Make room for the new "Definition of MSMEs" block.
#    The old "Sector Distribution Details" table (rows 24-43) needs to
#    shift down by 6 rows, to rows 30-49. Inserting 6 blank rows right
#    before the old row 24 achieves that and keeps every other cell's
#    value/style intact (including the hyperlink *cell* itself).
# ------------------------------------------------------------------
$ws.Rows("24:29").Insert()

# ------------------------------------------------------------------
# Helper: assign literal text, guarding against Excel's formula
# parser for values that start with '=', '+', '-' or '@'.
# ------------------------------------------------------------------
function Set-Text($rng, $text) {
    if ($text.Length -gt 0 -and ("=+-@".IndexOf($text.Substring(0,1)) -ge 0)) {
        $rng.Value = "'" + $text
    } else {
        $rng.Value = $text
    }
}

# ------------------------------------------------------------------
# 2. New table header (row 22) - bold "title" look.
# ------------------------------------------------------------------
Set-Text $ws.Range("B22") "Number of employees"
Set-Text $ws.Range("C22") "Assets (local currency, unless noted otherwise)"
Set-Text $ws.Range("D22") "Turnover (local currency, unless noted otherwise)"
$ws.Range("B22:D22").Font.Bold = $true

# ------------------------------------------------------------------
# 3. Data rows 23-26 - plain "Normal" look.
# ------------------------------------------------------------------
Set-Text $ws.Range("A23") "Micro"
Set-Text $ws.Range("B23") "=<19 Trade/Service"
$ws.Range("C23").Style = "Normal"
Set-Text $ws.Range("D23") "=< MNT 250 Millionlion Trade/Service"

Set-Text $ws.Range("A24") "Small"
Set-Text $ws.Range("B24") "=<19 Manufacturing<br/>=<49 Service"
$ws.Range("C24").Style = "Normal"
Set-Text $ws.Range("D24") "=< MNT 250 Millionlion Manufacturing<br/>=< MNT 1 Billionlion Service"

Set-Text $ws.Range("A25") "Medium"
Set-Text $ws.Range("B25") "=<149 Wholesale trade<br/>=<199 Retail Trade<br/>=<199 Manufacturing"
$ws.Range("C25").Style = "Normal"
Set-Text $ws.Range("D25") "=< MNT 1.5 Billionlion Wholesale trade<br/>=< MNT 1.5 Billionlion Retail Trade<br/>=< MNT 1.5 Billionlion Manufacturing"

Set-Text $ws.Range("A26") "Large"
Set-Text $ws.Range("B26") ">149 Wholesale trade<br/>>199 Retail Trade<br/>>199 Manufacturing"
$ws.Range("C26").Style = "Normal"
Set-Text $ws.Range("D26") "> MNT 1.5 Billionlion Wholesale trade<br/>> MNT 1.5 Billionlion Retail Trade<br/>> MNT 1.5 Billionlion Manufacturing"

$ws.Range("A23:D26").Font.Bold = $false
$ws.Range("A23:D26").Font.Italic = $false
$ws.Range("A23:D26").Font.Underline = $false

# ------------------------------------------------------------------
# 4. Fix the hyperlink: it used to live on the "source" citation cell
#    that has now moved from A38 to A44.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A44"), "http://www.mongolbank.mn/documents/moneypolicy/worldbank/developmentmodule/03.pdf")
$ws.Range("A44").Font.Underline = $true
$ws.Range("A44").Font.Color = 16711680
$ws.Range("A44").Font.Bold = $false
$ws.Range("A44").Font.Italic = $false
